$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new journal entries for rows 25-30 (2022-05-05, serial 44686) ---
# Row 25
$ws.Range("A25").Value = 44686
$ws.Range("B25").Value = 0.33333333333333331
$ws.Range("C25").Value = 0.37291666666666662

# Row 26
$ws.Range("A26").Value = 44686
$ws.Range("B26").Value = 0.37291666666666662
$ws.Range("C26").Value = 0.44236111111111115

# Row 27
$ws.Range("A27").Value = 44686
$ws.Range("B27").Value = 0.44236111111111115
$ws.Range("C27").Value = 0.4861111111111111

# Row 28
$ws.Range("A28").Value = 44686
$ws.Range("B28").Value = 0.4861111111111111
$ws.Range("C28").Value = 0.51041666666666663

# Row 30 (filled before row 29 so the shared-string table order matches the
# source workbook: the author apparently typed row 30 before going back to
# fill in row 29).
$ws.Range("A30").Value = 44686
$ws.Range("B30").Value = 0.62638888888888888
$ws.Range("C30").Value = 0.70486111111111116

# Descriptions (column E) / Solutions (column F)
$ws.Range("E25").Value = "Mise en forme final de la partie documentation Menu principale"
$ws.Range("E26").Value = "Documentation : Objectifs"
$ws.Range("E27").Value = "Documentation: Risque Technique"
$ws.Range("E28").Value = "Mise en forme Documentation"
$ws.Range("E30").Value = "Analyse de la partie Pathfinding"
$ws.Range("F30").Value = "Documentation sur le sujet afin de construire une explication brève concise et simple"

# Row 29 (filled last, after row 30, per the observed shared-string ordering)
$ws.Range("A29").Value = 44686
$ws.Range("B29").Value = 0.56388888888888888
$ws.Range("C29").Value = 0.62569444444444444
$ws.Range("E29").Value = "Remise en forme`ndu dossier de projet"

# --- View / selection changes ---
$ws.Range("A1:H30").Select()

# --- Page setup changes ---
$ws.PageSetup.LeftMargin = 18
$ws.PageSetup.RightMargin = 18
$ws.PageSetup.Zoom = 12
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
